$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 2112
$ws.Range("F5").Value = 7178
$ws.Range("F6").Value = 577
$ws.Range("F8").Value = 59
$ws.Range("F9").Value = 4684
$ws.Range("F10").Value = 6908
$ws.Range("F12").Value = 250
$ws.Range("F13").Value = 1448
$ws.Range("F14").Value = 836
$ws.Range("F15").Value = 138
$ws.Range("F16").Value = 38
$ws.Range("F17").Value = 1149
$ws.Range("F19").Value = 144
$ws.Range("F21").Value = 202
$ws.Range("F23").Value = 1110
$ws.Range("F24").Value = 748
$ws.Range("F25").Value = 43
$ws.Range("F26").Value = 1192
$ws.Range("F27").Value = 36
$ws.Range("F28").Value = 130
$ws.Range("F30").Value = 38
$ws.Range("F31").Value = 138
$ws.Range("F34").Value = 58
$ws.Range("F37").Value = 538
$ws.Range("F38").Value = 409
$ws.Range("F40").Value = 59
$ws.Range("F41").Value = 339
$ws.Range("F42").Value = 1193
$ws.Range("F43").Value = 556
$ws.Range("F44").Value = 75
$ws.Range("F45").Value = 131
$ws.Range("F46").Value = 14
$ws.Range("F47").Value = 13

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 24
$ws.Range("F12").Value = 24
$ws.Range("F16").Value = 1730
$ws.Range("F17").Value = 553
$ws.Range("F26").Value = 624
$ws.Range("F28").Value = 17
$ws.Range("F31").Value = 836
$ws.Range("F33").Value = 597
$ws.Range("F41").Value = 13

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 721
$ws.Range("F6").Value = 645
$ws.Range("F7").Value = 285
$ws.Range("F8").Value = 1474
$ws.Range("F9").Value = 2307

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 721
$ws.Range("F7").Value = 645
$ws.Range("F8").Value = 645
$ws.Range("F9").Value = 285
$ws.Range("F10").Value = 7178
$ws.Range("F11").Value = 59
$ws.Range("F12").Value = 4684
$ws.Range("F14").Value = 6908
$ws.Range("F15").Value = 250
$ws.Range("F16").Value = 1448
$ws.Range("F18").Value = 553
$ws.Range("F19").Value = 836
$ws.Range("F20").Value = 38
$ws.Range("F21").Value = 1149
$ws.Range("F22").Value = 144
$ws.Range("F23").Value = 1110
$ws.Range("F24").Value = 624
$ws.Range("F25").Value = 748
$ws.Range("F26").Value = 43
$ws.Range("F27").Value = 1192
$ws.Range("F28").Value = 17
$ws.Range("F30").Value = 836
$ws.Range("F32").Value = 58
$ws.Range("F35").Value = 538
$ws.Range("F36").Value = 597
$ws.Range("F37").Value = 409
$ws.Range("F39").Value = 59
$ws.Range("F41").Value = 339
$ws.Range("F42").Value = 556
$ws.Range("F45").Value = 13
$ws.Range("F46").Value = 131
$ws.Range("F48").Value = 14
$ws.Range("F49").Value = 13

